# Add some SMS examples to the "send_sms" form workbook.
#
# - survey sheet: insert a new "phone_number" text prompt before the
#   existing "send_sms" prompt row, and append two new prompt rows
#   ("odk_sms" and "odk_sms_automatic") after it.
# - prompt_types sheet: register the two new prompt types ("odk_sms",
#   "odk_sms_automatic") with type "integer" in the two already-blank
#   trailing rows.
# - finally leave the "prompt_types" sheet as the active tab/sheet.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$promptTypes = $wb.Worksheets.Item("prompt_types")

# --- survey sheet -----------------------------------------------------
# Existing layout (before):
#   Row1: type | name | display.text | display.hint | hideInContents
#   Row2: text | subject_name | Enter the subject's name.
#   Row3: send_sms | send_sms | This will send an sms to 360.430.1468. | <hint>
#
# Insert a blank row above the "send_sms" row first, so rows 5/6 below
# it already exist at their final positions, then fill row 5 (odk_sms),
# row 3 (phone_number) and row 6 (odk_sms_automatic) in that order, so
# the newly-created shared strings land in the same order as the target
# file.

$survey.Rows.Item(3).Insert()

$survey.Range("A5").Value = "odk_sms"
$survey.Range("B5").Value = "odk_sms"
$survey.Range("C5").Value = "This will send an sms via the sms bridge."

$survey.Range("A3").Value = "text"
$survey.Range("B3").Value = "phone_number"
$survey.Range("C3").Value = "Enter the phone number to which to send the text."

$survey.Range("A6").Value = "odk_sms_automatic"
$survey.Range("B6").Value = "odk_sms_automatic"
$survey.Range("C6").Value = "This will send an sms without requiring confirmation."

$survey.Range("A3:C3").Style = $survey.Range("A2:C2").Style
$survey.Rows.Item(3).RowHeight = 12
$survey.Rows.Item(5).RowHeight = 12.75
$survey.Rows.Item(6).RowHeight = 12.75

$survey.Columns.Item(2).ColumnWidth = 15.3

$survey.Range("C7").Select() | Out-Null

# --- prompt_types sheet ------------------------------------------------
# Existing layout (before):
#   Row1: prompt_type_name | type
#   Row2: send_sms | integer
#   Row3: (empty)
#   Row4: (empty)
#
# Fill the two already-blank trailing rows with the new prompt types.

$promptTypes.Range("A3").Value = "odk_sms"
$promptTypes.Range("B3").Value = "integer"

$promptTypes.Range("A4").Value = "odk_sms_automatic"
$promptTypes.Range("B4").Value = "integer"

$promptTypes.Range("B5").Select() | Out-Null

# --- activate the prompt_types tab -------------------------------------
$promptTypes.Activate()
